$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new value in A3, matching the target revision
$ws.Range("A3").Value = "test 2"

# Move the selection to the cell just after the new data, as in the original edit
$ws.Range("A4").Select()
